# Implemented getting number of lines for methods and classes.
$wb = $excel.ActiveWorkbook

# Add "classNumberOfLines" sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsClassLines = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsClassLines.Name = "classNumberOfLines"

$wsClassLines.Range("A1").Value = "Class Name"
$wsClassLines.Range("B1").Value = "Number of Lines"

$wsClassLines.Range("A2").Value = "com.macro.mall.config.SecuritySecureConfig"
$wsClassLines.Range("B2").NumberFormat = "@"
$wsClassLines.Range("B2").Value = "13"

$wsClassLines.Range("A3").Value = "com.macro.mall.MallMonitorApplication"
$wsClassLines.Range("B3").NumberFormat = "@"
$wsClassLines.Range("B3").Value = "6"

# Add "methodNumberOfLines" sheet right after "classNumberOfLines".
$wsMethodLines = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsClassLines)
$wsMethodLines.Name = "methodNumberOfLines"

$wsMethodLines.Range("A1").Value = "Class Name"
$wsMethodLines.Range("B1").Value = "Method Signature"
$wsMethodLines.Range("C1").Value = "Number of Lines"

$wsMethodLines.Range("A2").Value = "com.macro.mall.config.SecuritySecureConfig"
$wsMethodLines.Range("B2").Value = "configure(org.springframework.security.config.annotation.web.builders.HttpSecurity)"
$wsMethodLines.Range("C2").NumberFormat = "@"
$wsMethodLines.Range("C2").Value = "1"

$wsMethodLines.Range("A3").Value = "com.macro.mall.MallMonitorApplication"
$wsMethodLines.Range("B3").Value = "main(java.lang.String[])"
$wsMethodLines.Range("C3").NumberFormat = "@"
$wsMethodLines.Range("C3").Value = "3"
